# Daily attendance processing - reorder "Recorded By" (column G) names so that
# the literal token "System" (exact case) is moved to the end of the
# comma-separated list, preserving the relative order of every other item
# (including a lowercase "system" entry, which is left exactly where it was).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CharCodeEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) {
        return $false
    }
    for ($i = 0; $i -lt $s1.Length; $i++) {
        $c1 = [int][char]$s1[$i]
        $c2 = [int][char]$s2[$i]
        if ($c1 -ne $c2) {
            return $false
        }
    }
    return $true
}

function ReorderRecordedBy($s) {
    $parts = $s.Split(",")
    $trimmed = @()
    foreach ($p in $parts) {
        $trimmed += $p.Trim()
    }

    $nonSystem = @()
    $systemCount = 0
    foreach ($p in $trimmed) {
        if (CharCodeEquals $p "System") {
            $systemCount += 1
        } else {
            $nonSystem += $p
        }
    }

    for ($i = 0; $i -lt $systemCount; $i++) {
        $nonSystem += "System"
    }

    return ($nonSystem -join ", ")
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($current -eq $null) {
        continue
    }
    $current = "$current"
    if ($current -eq "") {
        continue
    }
    if ($current.IndexOf(",") -lt 0) {
        continue
    }

    $updated = ReorderRecordedBy $current
    if ($updated -ne $current) {
        $cell.Value2 = $updated
    }
}
